# Actualización de horarios - Línea 141 - 766
# Updates the "Última actualización" / "Total filas" headers plus refreshes
# the scrape timestamps/minute counts for already-present rows, and appends
# newly-scraped rows to each sheet.

$wb = $excel.ActiveWorkbook

$newTime = "03:57:17"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 20"

$ws1.Range("A16").Value = $newTime
$ws1.Range("D16").Value = 4

$ws1.Range("A18").Value = $newTime
$ws1.Range("D18").Value = 49

$ws1.Range("A19").Value = $newTime
$ws1.Range("D19").Value = 56

$ws1.Range("A21").Value = $newTime
$ws1.Range("D21").Value = 79

$ws1.Range("A22").Value = $newTime
$ws1.Range("B22").Value = "05:22"
$ws1.Range("C22").Value = "23_HERNANDEZ"
$ws1.Range("D22").Value = 85
$ws1.Range("E22").Value = "LP1912"

$ws1.Range("A23").Value = $newTime
$ws1.Range("B23").Value = "05:35"
$ws1.Range("C23").Value = "215B_EL PATO"
$ws1.Range("D23").Value = 98
$ws1.Range("E23").Value = "LP1912"

$ws1.Range("A24").Value = $newTime
$ws1.Range("B24").Value = "05:39"
$ws1.Range("C24").Value = "14_ABASTO"
$ws1.Range("D24").Value = 102
$ws1.Range("E24").Value = "LP1912"

$ws1.Range("A25").Value = $newTime
$ws1.Range("B25").Value = "05:46"
$ws1.Range("C25").Value = "15_ABASTO"
$ws1.Range("D25").Value = 109
$ws1.Range("E25").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 7"

$ws2.Range("A11").Value = $newTime
$ws2.Range("D11").Value = 49

$ws2.Range("A12").Value = $newTime
$ws2.Range("B12").Value = "05:35"
$ws2.Range("C12").Value = "215B_EL PATO"
$ws2.Range("D12").Value = 98
$ws2.Range("E12").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
$ws3.Range("A3").Value = "Total filas: 2"

$ws3.Range("A7").Value = $newTime
$ws3.Range("B7").Value = "05:44"
$ws3.Range("C7").Value = "215A_LA PLATA"
$ws3.Range("D7").Value = 107
$ws3.Range("E7").Value = "L6173"
